$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update review outcomes --------------------------------------------
# Rows 2,4,5,6,9,10 were reviewed & accepted -> Point status "Closed",
# Acceptance "Accepted". Rows 3,7,8 remain open / blank acceptance.
# Shared-string insertion order matters to reproduce the target sst index
# order (Closed=17, Accepted=18, comment=19), so write E (Closed) for all
# rows first, then D (Accepted) for all rows, then the new F3 comment.

$ws.Range("E2").Value = "Closed"
$ws.Range("E4").Value = "Closed"
$ws.Range("E5").Value = "Closed"
$ws.Range("E6").Value = "Closed"
$ws.Range("E9").Value = "Closed"
$ws.Range("E10").Value = "Closed"

$ws.Range("D2").Value = "Accepted"
$ws.Range("D4").Value = "Accepted"
$ws.Range("D5").Value = "Accepted"
$ws.Range("D6").Value = "Accepted"
$ws.Range("D9").Value = "Accepted"
$ws.Range("D10").Value = "Accepted"

# New reviewer comment on row 3 (F column). Match the left/top aligned,
# non-wrapping style already used by column C ("s=9") by copying the
# format from C4 before writing the text, so no new cellXfs are created.
$ws.Range("C4").Copy()
$ws.Range("F3").PasteSpecial(-4122)
$ws.Range("F3").Value = "Mali 13/3/2020: Point still open"
$excel.CutCopyMode = 0

# --- Re-apply the Open/Closed and Accepted/Rejected highlight rules ----
# so the workbook regenerates a fresh set of conditional-format dxf
# records (mirrors re-saving the review sheet's highlighting after the
# edits above).
$rE = $ws.Range("E2:E10")
$rD = $ws.Range("D2:D10")

for ($i = 1; $i -le 3; $i++) {
    $fc = $rE.FormatConditions.Add(9, , "Open")
    $fc.Font.Color = 6722871
    $fc.Interior.Color = 13561798

    $fc = $rE.FormatConditions.Add(9, , "Closed")
    $fc.Font.Color = 10284031

    $fc = $rE.FormatConditions.Add(9, , "Open")
    $fc.Font.Color = 10284031

    $fc = $rD.FormatConditions.Add(1, 3, '"Rejected"')
    $fc.Font.Color = 6722871
    $fc.Interior.Color = 13561798

    $fc = $rD.FormatConditions.Add(1, 3, '"Accepted"')
    $fc.Font.Color = 255
}

# --- Move the active selection (the review was scrolled back up to the
# top of the sheet and cell C3 selected once the new comment was added).
$ws.Range("C3").Select()
